# Insert a new weekly data row at row 6 (pushing existing rows 6-66 down to
# 7-67) and populate it with the new "Albahaca" price record, matching the
# target diff: dimension grows from A1:R66 to A1:R67, and all previously
# existing data rows 6-66 shift down by one row unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6; Excel shifts rows 6:66 -> 7:67 and copies the
# formatting (including the date style) from the row above, same as a normal
# interactive "Insert Row" in the UI.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value2 = 1
$ws.Range("B6").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value2 = "Arica y Parinacota"
$ws.Range("D6").Value2 = 45111
$ws.Range("E6").Value2 = 15
$ws.Range("F6").Value2 = 100112052
$ws.Range("G6").Value2 = "Albahaca"
$ws.Range("H6").Value2 = "Sin especificar"
$ws.Range("I6").Value2 = "Primera"
$ws.Range("J6").Value2 = 400
$ws.Range("K6").Value2 = 900
$ws.Range("L6").Value2 = 1000
$ws.Range("M6").Value2 = 962
$ws.Range("N6").Value2 = "$/paquete"
$ws.Range("O6").Value2 = "Región de Arica y Parinacota"
$ws.Range("P6").Value2 = 962
$ws.Range("Q6").Value2 = 1
$ws.Range("R6").Value2 = "Hortaliza"
